$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "Promo`rProbleemstelling`rDemo`rRealisaties`rTechnisch`rAanvullingen"
